# Update "想去人数" (want-to-go count) values that changed between the
# previous scrape and the new scrape (generated output at commit 456a3b4).
#
# Sheet "展览" (exhibitions) rows 4,5,14,16,18,20 column F
# Sheet "全部类型" (all types, aggregate) rows 4,12,36,38,41,43 column F
# "演出" and "本地生活" sheets are unaffected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1175
$ws1.Range("F5").Value = 18
$ws1.Range("F14").Value = 13184
$ws1.Range("F16").Value = 9
$ws1.Range("F18").Value = 5417
$ws1.Range("F20").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1175
$ws4.Range("F12").Value = 18
$ws4.Range("F36").Value = 13184
$ws4.Range("F38").Value = 9
$ws4.Range("F41").Value = 5417
$ws4.Range("F43").Value = 23
